$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text format first,
# so Excel keeps them as plain text strings (matching the source data which
# stores prices/percentages as text, not numbers).
$textCells = @("D5","D6","D8","D11","D14","D15","D19","D20","D21","D22","D23","D24","D26","D27","D30","D31","D33","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "91.818.34"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").Value = "3.163.58"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "239.93"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").Value = "621.72"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  +5.94%  "
$ws.Range("D8").Value = "0.377"
$ws.Range("E8").Value = "  +3.78%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "3.161.24"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("D11").Value = "0.745"
$ws.Range("E11").Value = "  +4.30%  "
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "35.61"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "5.59"
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("D16").Value = "91.275.40"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").Value = "3.752.62"
$ws.Range("D18").Value = "3.165.03"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").Value = "3.76"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "15.27"
$ws.Range("E20").Value = "  +10.95%  "
$ws.Range("D21").Value = "5.93"
$ws.Range("E21").Value = "  +10.14%  "
$ws.Range("D22").Value = "456.86"
$ws.Range("E22").Value = "  +5.73%  "
$ws.Range("D23").Value = "0.0000205"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("D24").Value = "9.22"
$ws.Range("E24").Value = "  +5.40%  "
$ws.Range("E25").Value = "  +8.49%  "
$ws.Range("D26").Value = "89.10"
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("D27").Value = "12.08"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "0.126"
$ws.Range("E30").Value = "  +40.52%  "
$ws.Range("D31").Value = "0.231"
$ws.Range("E31").Value = "  +17.97%  "
$ws.Range("E32").Value = "  +10.63%  "
$ws.Range("D33").Value = "9.43"
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("E34").Value = "  +14.47%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "0.932"
$ws.Range("E35").Value = "  -19.11%  "
$ws.Range("D36").Value = "7.70"
$ws.Range("E36").Value = "  +7.91%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "26.64"
$ws.Range("E37").Value = "  +4.19%  "
$ws.Range("D38").Value = "514.31"
$ws.Range("E38").Value = "  +3.85%  "
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").Value = "1.95"
$ws.Range("E39").Value = "  +3.65%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "1.36"
$ws.Range("E40").Value = "  +8.44%  "
$ws.Range("D41").Value = "0.454"
$ws.Range("E41").Value = "  +14.38%  "
$ws.Range("D42").Value = "3.83"
$ws.Range("E42").Value = "  +5.36%  "
$ws.Range("D43").Value = "3.48"
$ws.Range("E43").Value = "  -4.75%  "
$ws.Range("D44").Value = "22.20"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "0.715"
$ws.Range("E46").Value = "  +5.99%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "158.78"
$ws.Range("E47").Value = "  +3.95%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.95"
$ws.Range("E48").Value = "  +5.45%  "
$ws.Range("E49").Value = "  +6.22%  "
$ws.Range("E50").Value = "  +4.39%  "
$ws.Range("D51").Value = "44.12"
$ws.Range("E51").Value = "  -0.79%  "
